# Uncommented RAD Extension Payments Code and Test Data.
# The "Extension Payments" row (row 4) had its Execute flag set to
# "DONOTRUN"; re-enable it by setting the value back to "Y".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C4").Value = "Y"

# Reflect the new active cell selection on the sheet (matches the saved
# workbook state after making this edit).
$ws.Range("C4").Select()
